# Update the student list:
#  - LIU's password (D11) changes from "password" to "1" (stored as text)
#  - Trailing spaces are trimmed from several student names in column A
#  - The active selection moves to E12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set D11 to the literal text "1" (not a number) without altering cell
# formatting. We build it via a formula returning a text string, then
# paste back only the value (PasteSpecial xlPasteValues = -4163) so the
# cell keeps its original (default) style but becomes a text cell.
$cell = $ws.Range("D11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Trim trailing spaces from names
$ws.Range("A2").Value = "CHERN"
$ws.Range("A3").Value = "KOH"
$ws.Range("A4").Value = "BRANDON"
$ws.Range("A5").Value = "CALVIN"
$ws.Range("A10").Value = "LEE"

# Move the active selection to E12
$ws.Range("E12").Select()
